$d = $word.ActiveDocument

# Each of the three paragraphs below currently holds its sentence split
# across many single word/space runs (one <w:r> per word, plus one per
# separating space). This edit collapses each paragraph down to a single
# run holding the whole sentence, without changing the visible text or
# the paragraph style.
#
# Replacing a Range's .Text with a value that already matches the range's
# current (concatenated) text is treated as a no-op and would leave the
# old multi-run structure in place, so each paragraph is first stamped
# with a short placeholder to force a genuine content rewrite, and only
# then set to its real final text - that second write is what actually
# collapses the paragraph down to a single run.

$titleText = "Questions: Trigonometric identities (radians)"
$authorText = "Dzhemma Ruseva"
$abstractText = "A selection of questions on trigonometric identities, where angles are measured in radians."

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $styleName = $p.Style.NameLocal

    $finalText = $null
    if ($styleName -eq "Title") {
        $finalText = $titleText
    } elseif ($styleName -eq "Author") {
        $finalText = $authorText
    } elseif ($styleName -eq "Abstract") {
        $finalText = $abstractText
    }

    if ($finalText -ne $null) {
        $r = $p.Range
        $r.End = $r.End - 1
        $r.Text = "~"
        $r = $p.Range
        $r.End = $r.End - 1
        $r.Text = $finalText
    }
}
